$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Simple text substitutions (Find/Replace) ---

# Timestamp update
$d.Content.Find.Execute("Generated: 2025-11-21 18:26:37", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Generated: 2025-11-21 19:33:42", 2)

# Final scores list
$d.Content.Find.Execute("1. Alice: 100 points", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1. Mangals: 200 points", 2)

$d.Content.Find.Execute("2. Bob: -200 points", $true, $false, $false, $false, $false, `
    $true, 1, $false, "2. Badrie: -200 points", 2)

# Player headings (single run, no leading-space runs involved)
$d.Content.Find.Execute("Player: Alice", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Player: Mangals", 2)

$d.Content.Find.Execute("Player: Bob", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Player: Badrie", 2)

# --- New paragraph: "3. Arving: -300 points" after the Bob score line ---
$scoreLine = $d.Paragraphs(5)
$scoreLine.Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "3. Arving: -300 points"

# --- Replace the turn-by-turn detail paragraphs wholesale via InsertXML so
#     the leading-space runs keep their xml:space="preserve" marking. ---

# Mangals (formerly Alice) turn 1 detail -- paragraph 9 after the insert above
$mangalsTurnXml = '<w:p ' + $ns + '><w:r>' + `
    '<w:t>Turn 1:</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Category: Functions</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Question Value: 200 points</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Result: CORRECT</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Points Earned: +200</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Running Total: 200</w:t>' + `
    '</w:r></w:p>'
$d.Paragraphs(9).Range.InsertXML($mangalsTurnXml)

# Badrie (formerly Bob) turn 1 detail -- paragraph 11 after the insert above
$badrieTurnXml = '<w:p ' + $ns + '><w:r>' + `
    '<w:t>Turn 1:</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Category: Arrays</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Question Value: 200 points</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Result: INCORRECT</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Points Earned: -200</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Running Total: -200</w:t>' + `
    '</w:r></w:p>'
$d.Paragraphs(11).Range.InsertXML($badrieTurnXml)

# --- New paragraphs at the end: Player: Arving + turn breakdown ---
# Use a collapsed range just before the final paragraph mark so InsertXML
# appends a brand-new paragraph instead of overwriting an existing one.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$playerXml = '<w:p ' + $ns + '><w:r><w:rPr><w:b w:val="on"/></w:rPr><w:br/><w:t>Player: Arving</w:t></w:r></w:p>'
$endRange.InsertXML($playerXml)

$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$arvingTurnXml = '<w:p ' + $ns + '><w:r>' + `
    '<w:t>Turn 1:</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Category: Variables &amp; Data Types</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Question Value: 300 points</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Result: INCORRECT</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Points Earned: -300</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  Running Total: -300</w:t>' + `
    '</w:r></w:p>'
$endRange.InsertXML($arvingTurnXml)
